# Adds two new "Required Patients" rows (Patient 17 / Patient 18) to the
# test-suite grid, each describing the "deceased Datetime" access-period
# traits, replacing/extending the previous set of rows in Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column groups that share formatting (question-mark "?" cells), taken from
# the existing sheet so the new rows carry the same banded fill/border/font
# that the rest of the grid uses.
$style11Cols = @("B","E","G","H","K","M","O","Q")
$style12Cols = @("C","D","F","I","J","L","N","P")

# Prime the "Notes" column text for both new rows first (matches the shared
# -string insertion order of the authored workbook), then fill in the rest.
$ws.Range("R18").Copy()
$ws.Range("R20").PasteSpecial(-4122)
$ws.Range("R20").Value = "Patient Must have deceased Datetime with in allowed access period ""28"" days"

$ws.Range("R18").Copy()
$ws.Range("R19").PasteSpecial(-4122)
$ws.Range("R19").Value = "Patient Must have deceased Datetime over allowed access period ""28"" days. i.e deceased Datetime must be older than""28"" days."

# ---- Row 19: Patient 17 ----------------------------------------------
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A19").Value = "Patient 17"

foreach ($c in $style11Cols) {
    $ws.Range("B6").Copy()
    $ws.Range($c + "19").PasteSpecial(-4122)
    $ws.Range($c + "19").Value = "?"
}
foreach ($c in $style12Cols) {
    $ws.Range("F6").Copy()
    $ws.Range($c + "19").PasteSpecial(-4122)
    $ws.Range($c + "19").Value = "?"
}

$ws.Rows.Item(19).RowHeight = 60

# ---- Row 20: Patient 18 ----------------------------------------------
$ws.Range("A18").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A20").Value = "Patient 18"

foreach ($c in $style11Cols) {
    $ws.Range("B6").Copy()
    $ws.Range($c + "20").PasteSpecial(-4122)
    $ws.Range($c + "20").Value = "?"
}
foreach ($c in $style12Cols) {
    $ws.Range("F6").Copy()
    $ws.Range($c + "20").PasteSpecial(-4122)
    $ws.Range($c + "20").Value = "?"
}

$ws.Rows.Item(20).RowHeight = 30
